$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row: shift snapshot columns right, add Jun_15 / Jun_17 ---
# Row 1
$ws.Range('B1').Value = 'Jun_17'
$ws.Range('B1').ClearFormats()
$ws.Range('C1').Value = 'Jun_15'
$ws.Range('C1').ClearFormats()
$ws.Range('D1').Value = 'Jun_13'
$ws.Range('D1').ClearFormats()
$ws.Range('E1').Value = 'Jun_10'
$ws.Range('E1').ClearFormats()

# Row 2
$ws.Range('B2').Value = 'UN'
$ws.Range('B2').ClearFormats()
$ws.Range('C2').Value = 'UN'
$ws.Range('C2').ClearFormats()
$ws.Range('D2').Value = 'UN'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = 'UN'
$ws.Range('E2').ClearFormats()

# Row 3
$ws.Range('B3').Value = 'UN'
$ws.Range('B3').ClearFormats()
$ws.Range('C3').Value = 'UN'
$ws.Range('C3').ClearFormats()
$ws.Range('D3').Value = 'UN'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = 'UN'
$ws.Range('E3').ClearFormats()

# Row 4
$ws.Range('B4').Value = 'UN'
$ws.Range('B4').ClearFormats()
$ws.Range('C4').Value = 'UN'
$ws.Range('C4').ClearFormats()
$ws.Range('D4').Value = 'UN'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = 'UN'
$ws.Range('E4').ClearFormats()

# Row 5
$ws.Range('B5').Value = 'UN'
$ws.Range('B5').ClearFormats()
$ws.Range('C5').Value = 'UN'
$ws.Range('C5').ClearFormats()
$ws.Range('D5').Value = 'UN'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = 'UN'
$ws.Range('E5').ClearFormats()

# Row 6
$ws.Range('B6').Value = 'UN'
$ws.Range('B6').ClearFormats()
$ws.Range('C6').Value = 'UN'
$ws.Range('C6').ClearFormats()
$ws.Range('D6').Value = 'UN'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '10/11/2017,Reiterated Rating,Neutral,$25.00'
$ws.Range('E6').ClearFormats()

# Row 7
$ws.Range('B7').Value = 'UN'
$ws.Range('B7').ClearFormats()
$ws.Range('C7').Value = 'UN'
$ws.Range('C7').ClearFormats()
$ws.Range('D7').Value = 'UN'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = 'UN'
$ws.Range('E7').ClearFormats()

# Row 8
$ws.Range('B8').Value = 'UN'
$ws.Range('B8').ClearFormats()
$ws.Range('C8').Value = 'UN'
$ws.Range('C8').ClearFormats()
$ws.Range('D8').Value = 'UN'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = 'UN'
$ws.Range('E8').ClearFormats()

# Row 9
$ws.Range('B9').Value = 'UN'
$ws.Range('B9').ClearFormats()
$ws.Range('C9').Value = 'UN'
$ws.Range('C9').ClearFormats()
$ws.Range('D9').Value = 'UN'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = 'UN'
$ws.Range('E9').ClearFormats()

# Row 10
$ws.Range('B10').Value = 'UN'
$ws.Range('B10').ClearFormats()
$ws.Range('C10').Value = 'UN'
$ws.Range('C10').ClearFormats()
$ws.Range('D10').Value = 'UN'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '3/2/2018,Reiterated Rating,Overweight ➝ Overweight,$48.00'
$ws.Range('E10').ClearFormats()

# Row 11
$ws.Range('B11').Value = 'UN'
$ws.Range('B11').ClearFormats()
$ws.Range('C11').Value = 'UN'
$ws.Range('C11').ClearFormats()
$ws.Range('D11').Value = 'UN'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = 'UN'
$ws.Range('E11').ClearFormats()

# Row 12
$ws.Range('B12').Value = 'UN'
$ws.Range('B12').ClearFormats()
$ws.Range('C12').Value = 'UN'
$ws.Range('C12').ClearFormats()
$ws.Range('D12').Value = 'UN'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = 'UN'
$ws.Range('E12').ClearFormats()

# Row 13
$ws.Range('B13').Value = 'UN'
$ws.Range('B13').ClearFormats()
$ws.Range('C13').Value = 'UN'
$ws.Range('C13').ClearFormats()
$ws.Range('D13').Value = '6/11/2018,Upgrades,Hold -> Buy,$40.00 -> $50.00'
$ws.Range('D13').Interior.ColorIndex = 35
$ws.Range('E13').Value = '5/2/2018,Downgrades,Buy -> Hold,$21.00'
$ws.Range('E13').ClearFormats()

# Row 14
$ws.Range('B14').Value = 'UN'
$ws.Range('B14').ClearFormats()
$ws.Range('C14').Value = 'UN'
$ws.Range('C14').ClearFormats()
$ws.Range('D14').Value = 'UN'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = 'UN'
$ws.Range('E14').ClearFormats()

# Row 15
$ws.Range('B15').Value = 'UN'
$ws.Range('B15').ClearFormats()
$ws.Range('C15').Value = 'UN'
$ws.Range('C15').ClearFormats()
$ws.Range('D15').Value = 'UN'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = 'UN'
$ws.Range('E15').ClearFormats()

# Row 16
$ws.Range('B16').Value = 'UN'
$ws.Range('B16').ClearFormats()
$ws.Range('C16').Value = 'UN'
$ws.Range('C16').ClearFormats()
$ws.Range('D16').Value = 'UN'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = 'UN'
$ws.Range('E16').ClearFormats()

# Row 17
$ws.Range('B17').Value = 'UN'
$ws.Range('B17').ClearFormats()
$ws.Range('C17').Value = 'UN'
$ws.Range('C17').ClearFormats()
$ws.Range('D17').Value = 'UN'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '2/21/2018,Downgrade,Overweight ➝ Neutral,$42.00'
$ws.Range('E17').ClearFormats()

# Row 18
$ws.Range('B18').Value = 'UN'
$ws.Range('B18').ClearFormats()
$ws.Range('C18').Value = 'UN'
$ws.Range('C18').ClearFormats()
$ws.Range('D18').Value = 'UN'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = 'UN'
$ws.Range('E18').ClearFormats()

# Row 19
$ws.Range('B19').Value = 'UN'
$ws.Range('B19').ClearFormats()
$ws.Range('C19').Value = 'UN'
$ws.Range('C19').ClearFormats()
$ws.Range('D19').Value = 'UN'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '2/9/2018,Boost Price Target,Buy ➝ Buy,$35.00 ➝ $40.00'
$ws.Range('E19').Interior.ColorIndex = 35

# Row 20
$ws.Range('B20').Value = 'UN'
$ws.Range('B20').ClearFormats()
$ws.Range('C20').Value = 'UN'
$ws.Range('C20').ClearFormats()
$ws.Range('D20').Value = 'UN'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '5/2/2018,Lowers Target,Buy -> Buy,$54.00 -> $46.00'
$ws.Range('E20').ClearFormats()

# Row 21
$ws.Range('B21').Value = 'UN'
$ws.Range('B21').ClearFormats()
$ws.Range('C21').Value = 'UN'
$ws.Range('C21').ClearFormats()
$ws.Range('D21').Value = 'UN'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '3/9/2018,Downgrade,Outperform ➝ Market Perform,'
$ws.Range('E21').ClearFormats()

# Row 22
$ws.Range('B22').Value = 'UN'
$ws.Range('B22').ClearFormats()
$ws.Range('C22').Value = 'UN'
$ws.Range('C22').ClearFormats()
$ws.Range('D22').Value = 'UN'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = 'UN'
$ws.Range('E22').ClearFormats()

# Row 23
$ws.Range('B23').Value = 'UN'
$ws.Range('B23').ClearFormats()
$ws.Range('C23').Value = 'UN'
$ws.Range('C23').ClearFormats()
$ws.Range('D23').Value = 'UN'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = 'UN'
$ws.Range('E23').ClearFormats()

# Row 24
$ws.Range('B24').Value = 'UN'
$ws.Range('B24').ClearFormats()
$ws.Range('C24').Value = 'UN'
$ws.Range('C24').ClearFormats()
$ws.Range('D24').Value = 'UN'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = 'UN'
$ws.Range('E24').ClearFormats()

# Row 25
$ws.Range('B25').Value = 'UN'
$ws.Range('B25').ClearFormats()
$ws.Range('C25').Value = 'UN'
$ws.Range('C25').ClearFormats()
$ws.Range('D25').Value = 'UN'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = 'UN'
$ws.Range('E25').ClearFormats()

# Row 26
$ws.Range('B26').Value = 'UN'
$ws.Range('B26').ClearFormats()
$ws.Range('C26').Value = 'UN'
$ws.Range('C26').ClearFormats()
$ws.Range('D26').Value = 'UN'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = 'UN'
$ws.Range('E26').ClearFormats()

# Row 27
$ws.Range('B27').Value = 'UN'
$ws.Range('B27').ClearFormats()
$ws.Range('C27').Value = 'UN'
$ws.Range('C27').ClearFormats()
$ws.Range('D27').Value = 'UN'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '2/8/2018,Boost Price Target,Equal Weight ➝ Equal Weight,$30.00 ➝ $37.00'
$ws.Range('E27').Interior.ColorIndex = 35

# --- Column widths: C/D/E all width 8.0 (matches MarketBeat snapshot-column style) ---
$ws.Columns("C").ColumnWidth = 7.166666666666667
$ws.Columns("D").ColumnWidth = 7.166666666666667
$ws.Columns("E").ColumnWidth = 7.166666666666667
